# Change the template delimiters from "{ ... }" to "{% ... %}"
# (commit: "[delimiters][s] Delimiters changed to a more unique one.")
#
# The document body has 4 paragraphs:
#   1. (empty)
#   2. {#body loop}
#   3. {my body}
#   4. {/body loop}   (with a _GoBack bookmark split across two runs)
#
# The header and footer each have 3 (header) / 4 (footer) paragraphs with
# matching {#...}/{...}/{/...} tokens.
#
# We scope every Find/Execute to the specific paragraph's Range so that the
# short "{" / "}" / "body loop}" style fragments cannot cross-match text in
# a sibling paragraph.

$d = $word.ActiveDocument

function Replace-InRange($range, $old, $new) {
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# ---- Main document body ----
$paras = $d.Paragraphs

# Paragraph 2: {#body loop} -> {%#body loop%}
Replace-InRange $paras.Item(2).Range "{#body loop}" "{%#body loop%}"

# Paragraph 3: {my body} -> {%my body%}  (collapses the 3 runs into one)
Replace-InRange $paras.Item(3).Range "{my body}" "{%my body%}"

# Paragraph 4: {/ ... body loop} -> {%/ ... body loop%}
# Split in two finds so the bookmark in the middle of the paragraph is left untouched.
Replace-InRange $paras.Item(4).Range "{/" "{%/"
Replace-InRange $paras.Item(4).Range "body loop}" "body loop%}"

# ---- Headers / footers (every section) ----
foreach ($sec in $d.Sections) {
    $hdr = $sec.Headers.Item(1)
    $hParas = $hdr.Range.Paragraphs
    Replace-InRange $hParas.Item(1).Range "{#header loop}" "{%#header loop%}"
    Replace-InRange $hParas.Item(2).Range "{my header}" "{%my header%}"
    Replace-InRange $hParas.Item(3).Range "{/header loop}" "{%/header loop%}"

    $ftr = $sec.Footers.Item(1)
    $fParas = $ftr.Range.Paragraphs
    Replace-InRange $fParas.Item(1).Range "{#footer loop}" "{%#footer loop%}"
    Replace-InRange $fParas.Item(2).Range "{my footer}" "{%my footer%}"
    Replace-InRange $fParas.Item(3).Range "{/footer loop}" "{%/footer loop%}"
}

Write-Host "Body:" $d.Content.Text
Write-Host "Header:" $d.Sections.Item(1).Headers.Item(1).Range.Text
Write-Host "Footer:" $d.Sections.Item(1).Footers.Item(1).Range.Text
